$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: replace German terms with English translations (some with
# trailing space, per the source data).
$ws.Range("A2").Value = "Work "
$ws.Range("A3").Value = "Business"
$ws.Range("A4").Value = "Education"
$ws.Range("A5").Value = "Shopping "
$ws.Range("A6").Value = "Errands"
$ws.Range("A7").Value = "Leisure"
$ws.Range("A8").Value = "Accompaniment "

# Column D: new column holding the original German wording (row 1 repeats
# the "Wegzweck" header, row 2 mirrors the old "Quelle"-style bold cell).
$ws.Range("D1").Value = "Wegzweck"
$ws.Range("D2").Value = "Arbeit"
$ws.Range("D3").Value = "dienstlich"
$ws.Range("D4").Value = "Ausbildung"
$ws.Range("D5").Value = "Einkauf"
$ws.Range("D6").Value = "Erledigung"
$ws.Range("D7").Value = "Freizeit"
$ws.Range("D8").Value = "Begleitung"

# D2 picks up the same style C1 / old-A2 used to carry.
$ws.Range("D2").Font.Bold = $ws.Range("C1").Font.Bold

# The translated category cells wrap their (longer) text - new cell style.
$ws.Range("A2").WrapText = $true
$ws.Range("A5").WrapText = $true
$ws.Range("A8").WrapText = $true

# Leave the cursor where the author ended up after the edit.
$ws.Range("C13").Select() | Out-Null
